# Update "previsao_retorno.xlsx" (Resumo_por_Cliente sheet) with refreshed
# "Dados BIBI" snapshot: a handful of "situacao" (meses sem comprar) values
# were recomputed against a newer reference date, and two rows (54 and 112)
# had their purchase-count / date fields refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column J ("situacao") text updates -------------------------------
$ws.Range("J45").Value  = "INATIVO - 2.0 meses sem comprar"
$ws.Range("J73").Value  = "INATIVO - 32.9 meses sem comprar"
$ws.Range("J82").Value  = "INATIVO - 25.7 meses sem comprar"
$ws.Range("J87").Value  = "INATIVO - 4.4 meses sem comprar"
$ws.Range("J91").Value  = "INATIVO - 32.9 meses sem comprar"
$ws.Range("J92").Value  = "INATIVO - 13.2 meses sem comprar"
$ws.Range("J93").Value  = "INATIVO - 18.2 meses sem comprar"
$ws.Range("J98").Value  = "INATIVO - 1.7 meses sem comprar"
$ws.Range("J103").Value = "INATIVO - 24.4 meses sem comprar"
$ws.Range("J105").Value = "INATIVO - 24.8 meses sem comprar"
$ws.Range("J106").Value = "INATIVO - 15.2 meses sem comprar"
$ws.Range("J107").Value = "INATIVO - 6.1 meses sem comprar"
$ws.Range("J111").Value = "INATIVO - 8.0 meses sem comprar"

# --- Row 54: total_compras_historico / ultima_compra / proxima_compra --
$ws.Range("E54").Value = 28
$ws.Range("H54").Value = 45820.74543981482
$ws.Range("I54").Value = 45881.74543981482

# --- Row 112: total_compras_historico / ultima_compra / proxima_compra -
$ws.Range("E112").Value = 15572
$ws.Range("H112").Value = 45820.6597800926
$ws.Range("I112").Value = 45821.6597800926
